$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "...词级平均信息熵更大，而字符级平均..." -> "...词级平均信息熵和字符级平均..."
# ------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "词级平均信息熵更大，而字符级平均",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "词级平均信息熵和字符级平均", 2)
Write-Output "Edit1 Found: $found1"

# ------------------------------------------------------------------
# Edit 2: table value correction 0.0009386 -> 0.00009386
# (row "中文", column "平均信息熵（词）" of the first table)
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cell = $t.Cell(3, 3)
$found2 = $cell.Range.Find.Execute(
    "0.0009386",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "0.00009386", 2)
Write-Output "Edit2 Found: $found2"

# ------------------------------------------------------------------
# Edit 3a: "中文词级平均信息熵更高，反映了词汇组合的高度灵活性"
#          -> "中文词级平均信息熵更低，反映了词汇组合的高度灵活性"
# ------------------------------------------------------------------
$found3a = $d.Content.Find.Execute(
    "中文词级平均信息熵更高，反映了词汇组合的高度灵活性",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "中文词级平均信息熵更低，反映了词汇组合的高度灵活性", 2)
Write-Output "Edit3a Found: $found3a"

# ------------------------------------------------------------------
# Edit 3b: replace the rest of that sentence plus the following one
#          ("分词后...远超英文形态变化的词汇扩展模式。汉字单字多义的特性
#          在词级组合中引入歧义，需依赖上下文消除不确定性，从而推高熵值。")
#          with the new "中文词级...非线性叠加" sentence, and turn the
#          paragraph break that used to follow it into a manual line
#          break plus four spaces -- merging what used to be two
#          paragraphs into a single paragraph.
# ------------------------------------------------------------------
$zwnj = [char]0x200C
$oldQuoteOpen = [char]0x201C
$oldQuoteClose = [char]0x201D
$arrow = [char]0x2192

$old3b = "分词后词汇总量可达数万级（如" + $oldQuoteOpen + "人工智能" + $oldQuoteClose + $oldQuoteOpen + "区块链" + $oldQuoteClose + "等复合词），远超英文形态变化的词汇扩展模式" + $zwnj + "。汉字单字多义的特性在词级组合中引入歧义，需依赖上下文消除不确定性，从而推高熵值" + $zwnj + "。^p" + $zwnj + "中文字符与"
$new3b = "中文词级平均熵仍低于英文，但字符级与词级总信息熵更大，反映汉字系统通过多层级组合（单字" + $arrow + "复合词）实现信息密度的非线性叠加" + $zwnj + "^l    " + $zwnj + "中文字符与"

$found3b = $d.Content.Find.Execute(
    $old3b,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $new3b, 2)
Write-Output "Edit3b Found: $found3b"
